# [Fonds de solidarite] Add 2022-06-14 data
# Updates nombre_aides (col C) and montant_total (col E) for a set of rows
# on Sheet1 to reflect the latest data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new C (nombre_aides), new E (montant_total)
$updates = @(
    @{Row=7;   C=94343;  E=158181859}
    @{Row=8;   C=181364; E=652618365}
    @{Row=13;  C=37916;  E=97639732}
    @{Row=29;  C=77285;  E=428124618}
    @{Row=79;  C=14208;  E=57128900}
    @{Row=99;  C=136577; E=863158438}
    @{Row=125; C=4601;   E=13148080}
    @{Row=164; C=50581;  E=168933017}
    @{Row=168; C=285059; E=1211971202}
    @{Row=169; C=562632; E=1285123656}
    @{Row=170; C=367476; E=2847154599}
    @{Row=171; C=115191; E=447537564}
    @{Row=174; C=357297; E=1018855897}
    @{Row=175; C=125578; E=813925308}
    @{Row=179; C=235755; E=813149347}
    @{Row=210; C=6427;   E=19897530}
    @{Row=313; C=220660; E=1371177188}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C   # Column C = nombre_aides
    $ws.Cells.Item($u.Row, 5).Value = $u.E   # Column E = montant_total
}

$wb.Save()
